# Apply the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as *text* (e.g. "43.907.28", "0.1000") in the
# original workbook. A plain ".Value = ..." assignment lets Excel
# auto-coerce numeric-looking strings into real numbers (dropping things like
# trailing zeros), so we force the text format first, then reset the cell
# style back to "Normal" so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.907.28'
$ws.Range("E2").Value = '  -1.09%  '
Set-TextValue $ws.Range("D3") '2.347.95'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.28%  '
Set-TextValue $ws.Range("D5") '0.677'
$ws.Range("E5").Value = '  -0.04%  '
Set-TextValue $ws.Range("D6") '238.68'
$ws.Range("E6").Value = '  -0.04%  '
Set-TextValue $ws.Range("D7") '73.45'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.04%  '
Set-TextValue $ws.Range("D9") '0.592'
$ws.Range("E9").Value = '  +8.39%  '
Set-TextValue $ws.Range("D10") '0.1000'
$ws.Range("E10").Value = '  -2.98%  '
Set-TextValue $ws.Range("D11") '57.24'
$ws.Range("E11").Value = '  -0.24%  '
Set-TextValue $ws.Range("D12") '32.25'
$ws.Range("E12").Value = '  +9.56%  '
Set-TextValue $ws.Range("D13") '7.25'
$ws.Range("E13").Value = '  +7.96%  '
$ws.Range("E14").Value = '  +0.22%  '
Set-TextValue $ws.Range("D15") '2.697.71'
$ws.Range("E15").Value = '  -0.73%  '
Set-TextValue $ws.Range("D16") '16.50'
$ws.Range("E16").Value = '  -2.00%  '
Set-TextValue $ws.Range("D17") '0.894'
$ws.Range("E17").Value = '  -1.26%  '
Set-TextValue $ws.Range("D18") '2.344.67'
$ws.Range("E18").Value = '  -0.96%  '
Set-TextValue $ws.Range("D19") '43.830.24'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("E20").Value = '  -2.56%  '
Set-TextValue $ws.Range("D21") '6.72'
$ws.Range("E21").Value = '  +4.18%  '
Set-TextValue $ws.Range("D22") '76.67'
$ws.Range("E22").Value = '  -1.63%  '
Set-TextValue $ws.Range("D23") '258.40'
$ws.Range("E23").Value = '  +1.05%  '
$ws.Range("E24").Value = '  +22.16%  '
$ws.Range("E26").Value = '  -2.80%  '
$ws.Range("E27").Value = '  -2.09%  '
Set-TextValue $ws.Range("D28") '10.68'
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("E29").Value = '  +1.58%  '
Set-TextValue $ws.Range("D30") '22.55'
$ws.Range("E30").Value = '  +0.07%  '
Set-TextValue $ws.Range("D31") '175.40'
$ws.Range("E31").Value = '  +1.31%  '
Set-TextValue $ws.Range("D32") '0.128'
$ws.Range("E32").Value = '  -3.78%  '
$ws.Range("E33").Value = '  +2.18%  '
Set-TextValue $ws.Range("D34") '0.0759'
$ws.Range("E34").Value = '  +2.26%  '
Set-TextValue $ws.Range("D35") '5.19'
$ws.Range("E35").Value = '  -0.39%  '
Set-TextValue $ws.Range("D36") '5.49'
$ws.Range("E36").Value = '  +5.23%  '
Set-TextValue $ws.Range("D37") '3.73'
$ws.Range("E37").Value = '  -5.19%  '
Set-TextValue $ws.Range("D38") '2.34'
$ws.Range("E38").Value = '  -4.04%  '
Set-TextValue $ws.Range("D39") '6.27'
$ws.Range("E39").Value = '  -3.54%  '
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("E41").Value = '  +12.37%  '
Set-TextValue $ws.Range("D42") '0.204'
$ws.Range("E42").Value = '  +11.70%  '
Set-TextValue $ws.Range("D43") '18.82'
$ws.Range("E43").Value = '  -4.72%  '
$ws.Range("B44").Value = 'BinanceUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D44") '1.00'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D45") '8.93'
$ws.Range("E45").Value = '  +0.41%  '
Set-TextValue $ws.Range("D46") '4.69'
$ws.Range("E46").Value = '  +4.34%  '
Set-TextValue $ws.Range("D47") '2.50'
$ws.Range("E47").Value = '  +6.10%  '
Set-TextValue $ws.Range("D48") '57.64'
$ws.Range("E48").Value = '  +9.19%  '
Set-TextValue $ws.Range("D49") '1.23'
$ws.Range("E49").Value = '  -1.99%  '
Set-TextValue $ws.Range("D50") '1.16'
$ws.Range("E50").Value = '  -0.38%  '
Set-TextValue $ws.Range("D51") '99.45'
$ws.Range("E51").Value = '  +0.79%  '
